$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D2").Value = 0.4191875549243985
$ws.Range("E2").Value = 0.1988962653115924
$ws.Range("F2").Value = 0.1513854117655921
$ws.Range("H2").Value = 0.02876653646983837
$ws.Range("I2").Value = 0.006646683062633787
$ws.Range("L2").Value = 0.02304694628723282
$ws.Range("M2").Value = 0.03217939539793669
$ws.Range("O2").Value = 0.05221316467841435
$ws.Range("P2").Value = 0.01936879692521693
$ws.Range("Q2").Value = 0.02897323996792089
$ws.Range("R2").Value = 0.01872522003150782
$ws.Range("X2").Value = 0.006331095914911082
$ws.Range("AA2").Value = 0.01015584799425352
$ws.Range("AC2").Value = 0.003888675663626032
$ws.Range("AD2").Value = 0.0002351656049246838
$ws.Range("G3").Value = 0.03105458991453655
$ws.Range("H3").Value = 0.1678595640736542
$ws.Range("I3").Value = 0.34417914417735
$ws.Range("J3").Value = 0.104096976904886
$ws.Range("K3").Value = 0.002399835875633097
$ws.Range("L3").Value = 0.0343734766453744
$ws.Range("N3").Value = 0.0203566103724366
$ws.Range("P3").Value = 0.02836404672136914
$ws.Range("Q3").Value = 0.04088274584005891
$ws.Range("S3").Value = 0.1256742090586891
$ws.Range("U3").Value = 0.07812054548280759
$ws.Range("W3").Value = 0.02263825493320439
$ws.Range("E4").Value = 0.2102933323559729
$ws.Range("G4").Value = 0.2551194418598071
$ws.Range("I4").Value = 0.01196426799722883
$ws.Range("J4").Value = 0.06091237193315317
$ws.Range("N4").Value = 0.02476182953853135
$ws.Range("O4").Value = 0.03527098174367467
$ws.Range("P4").Value = 0.016806218984216
$ws.Range("Q4").Value = 0.07741725876040371
$ws.Range("S4").Value = 0.1814439198682341
$ws.Range("T4").Value = 0.02304528828274621
$ws.Range("U4").Value = 0.07985114096158373
$ws.Range("V4").Value = 0.01644086893556386
$ws.Range("W4").Value = 0.0001924261553083983
$ws.Range("Z4").Value = 0.003530742484204563
$ws.Range("AF4").Value = 0.002651252489972503
$ws.Range("AG4").Value = 0.0002986576493988193
$ws.Range("E5").Value = 0.038575329251885
$ws.Range("F5").Value = 0.3105201763017616
$ws.Range("G5").Value = 0.2347095428401411
$ws.Range("H5").Value = 0.0535557580681697
$ws.Range("J5").Value = 0.04135272014704477
$ws.Range("K5").Value = 0.01258300349277528
$ws.Range("L5").Value = 0.04308974064189475
$ws.Range("M5").Value = 0.003722742016594229
$ws.Range("N5").Value = 0.01525273473164749
$ws.Range("O5").Value = 0.03983807078402234
$ws.Range("Q5").Value = 0.07343847892181843
$ws.Range("S5").Value = 0.1246601080673956
$ws.Range("U5").Value = 0.003025394034302868
$ws.Range("Y5").Value = 0.003381388787957729
$ws.Range("AI5").Value = 0.002294811912589189
$ws.Range("D6").Value = 0.3936349018295414
$ws.Range("E6").Value = 0.05602963586919978
$ws.Range("F6").Value = 0.169089399176386
$ws.Range("H6").Value = 0.02309635216080352
$ws.Range("I6").Value = 0.0288561171074857
$ws.Range("M6").Value = 0.07281874308068702
$ws.Range("O6").Value = 0.03022702261641498
$ws.Range("P6").Value = 0.03291564040568214
$ws.Range("Q6").Value = 0.1207201342842007
$ws.Range("R6").Value = 0.03603122104201079
$ws.Range("S6").Value = 0.03647625380050649
$ws.Range("U6").Value = 0.0001045786270813876

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D2").Value = 0.4191875549243985
$ws.Range("E2").Value = 0.618083820235991
$ws.Range("F2").Value = 0.7694692320015831
$ws.Range("G2").Value = 0.7694692320015831
$ws.Range("H2").Value = 0.7982357684714215
$ws.Range("I2").Value = 0.8048824515340554
$ws.Range("J2").Value = 0.8048824515340554
$ws.Range("K2").Value = 0.8048824515340554
$ws.Range("L2").Value = 0.8279293978212882
$ws.Range("M2").Value = 0.8601087932192248
$ws.Range("N2").Value = 0.8601087932192248
$ws.Range("O2").Value = 0.9123219578976391
$ws.Range("P2").Value = 0.9316907548228561
$ws.Range("Q2").Value = 0.960663994790777
$ws.Range("R2").Value = 0.9793892148222848
$ws.Range("S2").Value = 0.9793892148222848
$ws.Range("T2").Value = 0.9793892148222848
$ws.Range("U2").Value = 0.9793892148222848
$ws.Range("V2").Value = 0.9793892148222848
$ws.Range("W2").Value = 0.9793892148222848
$ws.Range("X2").Value = 0.9857203107371959
$ws.Range("Y2").Value = 0.9857203107371959
$ws.Range("Z2").Value = 0.9857203107371959
$ws.Range("AA2").Value = 0.9958761587314494
$ws.Range("AB2").Value = 0.9958761587314494
$ws.Range("AC2").Value = 0.9997648343950755
$ws.Range("AD2").Value = 1
$ws.Range("AE2").Value = 1
$ws.Range("AF2").Value = 1
$ws.Range("AG2").Value = 1
$ws.Range("AH2").Value = 1
$ws.Range("AI2").Value = 1
$ws.Range("AJ2").Value = 1
$ws.Range("AK2").Value = 1
$ws.Range("AL2").Value = 1
$ws.Range("G3").Value = 0.03105458991453655
$ws.Range("H3").Value = 0.1989141539881908
$ws.Range("I3").Value = 0.5430932981655407
$ws.Range("J3").Value = 0.6471902750704268
$ws.Range("K3").Value = 0.6495901109460599
$ws.Range("L3").Value = 0.6839635875914343
$ws.Range("M3").Value = 0.6839635875914343
$ws.Range("N3").Value = 0.7043201979638709
$ws.Range("O3").Value = 0.7043201979638709
$ws.Range("P3").Value = 0.7326842446852401
$ws.Range("Q3").Value = 0.773566990525299
$ws.Range("R3").Value = 0.773566990525299
$ws.Range("S3").Value = 0.8992411995839882
$ws.Range("T3").Value = 0.8992411995839882
$ws.Range("U3").Value = 0.9773617450667957
$ws.Range("V3").Value = 0.9773617450667957
$ws.Range("E4").Value = 0.2102933323559729
$ws.Range("F4").Value = 0.2102933323559729
$ws.Range("G4").Value = 0.46541277421578
$ws.Range("H4").Value = 0.46541277421578
$ws.Range("I4").Value = 0.4773770422130088
$ws.Range("J4").Value = 0.5382894141461619
$ws.Range("K4").Value = 0.5382894141461619
$ws.Range("L4").Value = 0.5382894141461619
$ws.Range("M4").Value = 0.5382894141461619
$ws.Range("N4").Value = 0.5630512436846933
$ws.Range("O4").Value = 0.598322225428368
$ws.Range("P4").Value = 0.615128444412584
$ws.Range("Q4").Value = 0.6925457031729877
$ws.Range("R4").Value = 0.6925457031729877
$ws.Range("S4").Value = 0.8739896230412217
$ws.Range("T4").Value = 0.897034911323968
$ws.Range("U4").Value = 0.9768860522855517
$ws.Range("V4").Value = 0.9933269212211155
$ws.Range("W4").Value = 0.9935193473764239
$ws.Range("X4").Value = 0.9935193473764239
$ws.Range("Y4").Value = 0.9935193473764239
$ws.Range("Z4").Value = 0.9970500898606285
$ws.Range("AA4").Value = 0.9970500898606285
$ws.Range("AB4").Value = 0.9970500898606285
$ws.Range("AC4").Value = 0.9970500898606285
$ws.Range("AD4").Value = 0.9970500898606285
$ws.Range("AE4").Value = 0.9970500898606285
$ws.Range("AF4").Value = 0.999701342350601
$ws.Range("AG4").Value = 0.9999999999999999
$ws.Range("AH4").Value = 0.9999999999999999
$ws.Range("AI4").Value = 0.9999999999999999
$ws.Range("AJ4").Value = 0.9999999999999999
$ws.Range("AK4").Value = 0.9999999999999999
$ws.Range("AL4").Value = 0.9999999999999999
$ws.Range("E5").Value = 0.038575329251885
$ws.Range("F5").Value = 0.3490955055536466
$ws.Range("G5").Value = 0.5838050483937877
$ws.Range("H5").Value = 0.6373608064619574
$ws.Range("I5").Value = 0.6373608064619574
$ws.Range("J5").Value = 0.6787135266090021
$ws.Range("K5").Value = 0.6912965301017774
$ws.Range("L5").Value = 0.7343862707436721
$ws.Range("M5").Value = 0.7381090127602664
$ws.Range("N5").Value = 0.7533617474919139
$ws.Range("O5").Value = 0.7931998182759362
$ws.Range("P5").Value = 0.7931998182759362
$ws.Range("Q5").Value = 0.8666382971977546
$ws.Range("R5").Value = 0.8666382971977546
$ws.Range("S5").Value = 0.9912984052651502
$ws.Range("T5").Value = 0.9912984052651502
$ws.Range("U5").Value = 0.9943237992994531
$ws.Range("V5").Value = 0.9943237992994531
$ws.Range("W5").Value = 0.9943237992994531
$ws.Range("X5").Value = 0.9943237992994531
$ws.Range("Y5").Value = 0.9977051880874108
$ws.Range("Z5").Value = 0.9977051880874108
$ws.Range("AA5").Value = 0.9977051880874108
$ws.Range("AB5").Value = 0.9977051880874108
$ws.Range("AC5").Value = 0.9977051880874108
$ws.Range("AD5").Value = 0.9977051880874108
$ws.Range("AE5").Value = 0.9977051880874108
$ws.Range("AF5").Value = 0.9977051880874108
$ws.Range("AG5").Value = 0.9977051880874108
$ws.Range("AH5").Value = 0.9977051880874108
$ws.Range("D6").Value = 0.3936349018295414
$ws.Range("E6").Value = 0.4496645376987412
$ws.Range("F6").Value = 0.6187539368751271
$ws.Range("G6").Value = 0.6187539368751271
$ws.Range("H6").Value = 0.6418502890359307
$ws.Range("I6").Value = 0.6707064061434164
$ws.Range("J6").Value = 0.6707064061434164
$ws.Range("K6").Value = 0.6707064061434164
$ws.Range("L6").Value = 0.6707064061434164
$ws.Range("M6").Value = 0.7435251492241034
$ws.Range("N6").Value = 0.7435251492241034
$ws.Range("O6").Value = 0.7737521718405184
$ws.Range("P6").Value = 0.8066678122462005
$ws.Range("Q6").Value = 0.9273879465304012
$ws.Range("R6").Value = 0.9634191675724121
$ws.Range("S6").Value = 0.9998954213729185
$ws.Range("T6").Value = 0.9998954213729185
$ws.Range("U6").Value = 0.9999999999999999
$ws.Range("V6").Value = 0.9999999999999999
$ws.Range("W6").Value = 0.9999999999999999
$ws.Range("X6").Value = 0.9999999999999999
$ws.Range("Y6").Value = 0.9999999999999999
$ws.Range("Z6").Value = 0.9999999999999999
$ws.Range("AA6").Value = 0.9999999999999999
$ws.Range("AB6").Value = 0.9999999999999999
$ws.Range("AC6").Value = 0.9999999999999999
$ws.Range("AD6").Value = 0.9999999999999999
$ws.Range("AE6").Value = 0.9999999999999999
$ws.Range("AF6").Value = 0.9999999999999999
$ws.Range("AG6").Value = 0.9999999999999999
$ws.Range("AH6").Value = 0.9999999999999999
$ws.Range("AI6").Value = 0.9999999999999999
$ws.Range("AJ6").Value = 0.9999999999999999
$ws.Range("AK6").Value = 0.9999999999999999
$ws.Range("AL6").Value = 0.9999999999999999

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("D2").Value = 4
$ws.Range("F2").Value = 0.618083820235991
$ws.Range("G2").Value = 3
$ws.Range("F3").Value = 0.5430932981655407
$ws.Range("D4").Value = 9
$ws.Range("F4").Value = 0.5382894141461619
$ws.Range("G4").Value = 7
$ws.Range("F5").Value = 0.5838050483937877
$ws.Range("D6").Value = 5
$ws.Range("F6").Value = 0.6187539368751271
$ws.Range("G6").Value = 4

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D2").Value = 5
$ws.Range("F2").Value = 0.7694692320015831
$ws.Range("G2").Value = 4
$ws.Range("D3").Value = 13
$ws.Range("F3").Value = 0.7043201979638709
$ws.Range("G3").Value = 8
$ws.Range("D4").Value = 18
$ws.Range("F4").Value = 0.8739896230412217
$ws.Range("G4").Value = 16
$ws.Range("D5").Value = 11
$ws.Range("F5").Value = 0.7343862707436721
$ws.Range("G5").Value = 8
$ws.Range("D6").Value = 12
$ws.Range("F6").Value = 0.7435251492241034
$ws.Range("G6").Value = 11

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D2").Value = 8
$ws.Range("F2").Value = 0.8048824515340554
$ws.Range("G2").Value = 7
$ws.Range("F3").Value = 0.8992411995839882
$ws.Range("F4").Value = 0.8739896230412217
$ws.Range("F5").Value = 0.8666382971977546
$ws.Range("D6").Value = 15
$ws.Range("F6").Value = 0.8066678122462005
$ws.Range("G6").Value = 14

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("D2").Value = 14
$ws.Range("F2").Value = 0.9123219578976391
$ws.Range("G2").Value = 13
$ws.Range("D3").Value = 20
$ws.Range("F3").Value = 0.9773617450667957
$ws.Range("G3").Value = 15
$ws.Range("D4").Value = 20
$ws.Range("F4").Value = 0.9768860522855517
$ws.Range("G4").Value = 18
$ws.Range("F5").Value = 0.9912984052651502
$ws.Range("F6").Value = 0.9273879465304012
